# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.046.33"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "'2.312.41"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'304.09"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'100.97"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("E9").Value = "  +3.71%  "
$ws.Range("D10").Value = "'34.94"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").Value = "'17.96"
$ws.Range("E13").Value = "  +15.14%  "
$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "'2.689.04"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "'2.363.01"
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("E17").Value = "  +4.56%  "
$ws.Range("D18").Value = "'42.980.83"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").Value = "'12.50"
$ws.Range("E19").Value = "  +7.80%  "
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "'0.0₃0906"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "'67.76"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "'237.42"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  +12.56%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'24.87"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.10"
$ws.Range("E28").Value = "  -8.53%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'168.01"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'34.12"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'9.18"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("D36").Value = "'17.10"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'0.0692"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").Value = "'1.80"
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("D43").Value = "'2.004.27"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").Value = "'10.19"
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("D46").Value = "'17.54"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").Value = "'55.77"
$ws.Range("E48").Value = "  +7.48%  "
$ws.Range("D49").Value = "'2.531.44"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +5.56%  "
$ws.Range("E51").Value = "  +0.96%  "
